$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above row 441; the existing data previously on
# rows 441:504 shifts down to 444:507 (all other rows/cols untouched).
$ws.Range("441:443").Insert()

# Row 444 now holds the template that used to live on row 441 (same
# Mercado/Region/Producto/Categoria/Variedad/Unidad/Origen/Kg-unidad
# columns for every row in this block) - clone it into the 3 fresh rows
# so every non-varying column (A,B,C,E:L,Q,R,T) is populated correctly,
# then overwrite the per-row values that actually differ (D,L,M,N,O,P,S).
$ws.Range("A444:T444").Copy()
$ws.Range("A441:T441").PasteSpecial()
$ws.Range("A444:T444").Copy()
$ws.Range("A442:T442").PasteSpecial()
$ws.Range("A444:T444").Copy()
$ws.Range("A443:T443").PasteSpecial()

# New week of Kiwi price data (fecha serial 44522) for Femacal de La Calera.
# Row 441: Especial
$ws.Range("D441").Value = 44522
$ws.Range("L441").Value = "Especial"
$ws.Range("M441").Value = 65
$ws.Range("N441").Value = 12000
$ws.Range("O441").Value = 12000
$ws.Range("P441").Value = 12000
$ws.Range("S441").Value = 1200

# Row 442: Primera
$ws.Range("D442").Value = 44522
$ws.Range("L442").Value = "Primera"
$ws.Range("M442").Value = 60
$ws.Range("N442").Value = 11000
$ws.Range("O442").Value = 11000
$ws.Range("P442").Value = 11000
$ws.Range("S442").Value = 1100

# Row 443: Segunda
$ws.Range("D443").Value = 44522
$ws.Range("L443").Value = "Segunda"
$ws.Range("M443").Value = 65
$ws.Range("N443").Value = 9000
$ws.Range("O443").Value = 9000
$ws.Range("P443").Value = 9000
$ws.Range("S443").Value = 900
